$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update Curvature ("abs_curv", column I) values using JM's methods ---
$ws.Range("I5").Value  = 3.3065552899999999
$ws.Range("I6").Value  = 5.3897394939999996
$ws.Range("I7").Value  = 16.494041790000001
$ws.Range("I11").Value = 3.951547315
$ws.Range("I12").Value = 2.658635834
$ws.Range("I13").Value = 11.362216119999999

# --- Re-style those same cells onto a shared "Arial 11" look (no explicit
#     number format override), matching the single new cellXfs entry that
#     replaces the two old ones. Reset to the base style first so the old
#     number format / alignment don't linger, then apply the new font. ---
$cells = @("I5", "I6", "I7", "I11", "I12", "I13")
foreach ($cell in $cells) {
    $rng = $ws.Range($cell)
    $rng.Style = "Normal"
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 11
}

# --- View state: scroll/selection moved to J21 ---
$ws.Activate()
$ws.Range("J21").Select()
